# Adicionados balancos concatenados em uma unica planilha.
# Adds columns AO (31/03/2024) and AP (30/06/2024) to Sheet1,
# mirroring the existing quarterly-balance layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new quarter labels, matching AN1 header style ---
$ws.Range("AN1").Copy()
$ws.Range("AO1").PasteSpecial(-4122)
$ws.Range("AO1").Value2 = "31/03/2024"

$ws.Range("AN1").Copy()
$ws.Range("AP1").PasteSpecial(-4122)
$ws.Range("AP1").Value2 = "30/06/2024"

# --- Data rows: (row, AO value, AP value) ---
$data = @(
    ,(2, 1291970.048, 1306680.96)
    ,(3, 832136, 837590.976)
    ,(4, 349251.008, 343246.016)
    ,(5, 0, 0)
    ,(6, 202911.008, 198092.992)
    ,(7, 222414, 249752)
    ,(8, 0, 0)
    ,(9, 34190, 37667)
    ,(10, 0, 0)
    ,(11, 23370, 8833)
    ,(12, 23807, 32644)
    ,(13, 0, 0)
    ,(14, 0, 0)
    ,(15, 0, 0)
    ,(16, 0, 0)
    ,(17, 0, 0)
    ,(18, 0, 0)
    ,(19, 21121, 30138)
    ,(20, 0, 0)
    ,(21, 0, 0)
    ,(22, 0, 0)
    ,(23, 339439.008, 336262.016)
    ,(24, 96588, 100184)
    ,(25, 0, 0)
    ,(26, 1291970.048, 1306680.96)
    ,(27, 239284, 271737.984)
    ,(28, 33502, 38749)
    ,(29, 82003, 104517)
    ,(30, 8359, 13941)
    ,(31, 98208, 96511)
    ,(32, 263, 334)
    ,(33, 0, 0)
    ,(34, 16949, 17686)
    ,(35, 0, 0)
    ,(36, 0, 0)
    ,(37, 349872, 345459.008)
    ,(38, 323216, 314280.992)
    ,(39, 0, 0)
    ,(40, 21708, 23789)
    ,(41, 0, 0)
    ,(42, 0, 0)
    ,(43, 4948, 7389)
    ,(44, 0, 0)
    ,(45, 0, 0)
    ,(46, 20, 19)
    ,(47, 702794.0159999999, 689465.032)
    ,(48, 599822.976, 599822.976)
    ,(49, 1998, 2294)
    ,(50, 0, 0)
    ,(51, 69708, 38708)
    ,(52, 12852, 29748)
    ,(53, 18413, 18892)
    ,(54, 0, 0)
    ,(55, 0, 0)
    ,(56, 0, 0)
    ,(59, 178390, 217214)
    ,(60, -92796, -112605)
    ,(61, 85594, 104609)
    ,(62, -57070, -64907)
    ,(63, -13224, -13869)
    ,(64, 5550, 0)
    ,(65, 6452, 3529)
    ,(66, -902, -3102)
    ,(67, 0, 0)
    ,(68, -1236, -1277)
    ,(69, 8775, 13458)
    ,(70, -10011, -14735)
    ,(74, 25164, 24983)
    ,(75, -5963, -17060)
    ,(76, -802, 8973)
    ,(79, 3, 0)
    ,(80, 12852, 16896)
)

foreach ($item in $data) {
    $r = $item[0]
    $aoVal = $item[1]
    $apVal = $item[2]
    $ws.Cells.Item($r, 41).Value2 = $aoVal
    $ws.Cells.Item($r, 42).Value2 = $apVal
}

